$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 3 (shifts existing rows 3-12 down to 4-13) ---
$ws.Rows.Item(3).Insert()

# --- Copy number formatting (style) from row 2's B/C cells onto the new row 3 cells ---
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Re-create the fill-down formulas for the shifted rows 4:13 (previously rows 3:12),
#     one row at a time so the formula engine keeps each row's cached value correct ---
for ($r = 4; $r -le 13; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=A$r"
    $ws.Cells.Item($r, 6).Formula = "=(E$r*D$r*SQRT(3))/1000000"
}

# --- Populate new row 3 with data for the 419 line type ---
$ws.Range("A3").Value = 419
$ws.Range("B3").Value = 0.13
$ws.Range("C3").Value = 0.638
$ws.Range("D3").Value = 419
$ws.Range("E3").Value = 21000
$ws.Range("F3").Formula = "=(E3*D3*SQRT(3))/1000000"

# --- Update scaled values for row 2 (SBase rescale of r/x parameters) ---
$ws.Range("B2").Value = 0.084994033412887832
$ws.Range("C2").Value = 0.41414671814671816

# --- Update the active selection shown when the workbook is opened ---
$ws.Range("J6").Select()
